$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new row of data (row 13) following the same pattern as the
# preceding rows: a date in column A and daily expense figures in B:M.
# Copy A12's formatting (date number format) down into A13, then set
# the new row's date value.
$ws.Range("A12").Copy($ws.Range("A13"))
$ws.Range("A13").Value = 43801

$ws.Range("B13").Value = 0
$ws.Range("C13").Value = 0
$ws.Range("D13").Value = 0
$ws.Range("E13").Value = 0
$ws.Range("F13").Value = 0
$ws.Range("G13").Value = 0
$ws.Range("H13").Value = 0
$ws.Range("I13").Value = 0
$ws.Range("J13").Value = 955
$ws.Range("K13").Value = 1
$ws.Range("L13").Value = 2
$ws.Range("M13").Value = 3

# Match the selection state recorded after the edit.
$ws.Range("J13").Select()
